$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25, columns B,C,D,E,F,I,J,K,L,M,N
# (column G stays 1, column H has no data, column A is the bus index - unchanged)
$newValues = @{
    "2" = @(1.02, 1.015994122660331, 1.04600280854966, 1.017522245995445, 1.0478698538534, 1.03560001447132, 1.021216464669534, 1.048769464460674, 1.0203711857041, 1.050631285144986, 1.01100072096038)
    "3" = @(1.02, 1.017528445337316, 1.04676595276577, 1.018843091799321, 1.048965554791632, 1.035749323793917, 1.022382965011307, 1.049344289238011, 1.021496222811384, 1.051538182075071, 1.011403083142908)
    "4" = @(1.02, 1.018520050010412, 1.047258516485378, 1.019696946306224, 1.049673394045788, 1.035844112666199, 1.02313626962621, 1.049714339045488, 1.022222851166886, 1.052123264913848, 1.011662484048454)
    "5" = @(1.02, 1.018936640393844, 1.047465293424365, 1.020055716127835, 1.049970695881091, 1.035883525271962, 1.02345260576678, 1.049869453636641, 1.022528008928907, 1.05236881908561, 1.011771309551694)
    "6" = @(1.02, 1.019006571557109, 1.047499994786048, 1.020115944154147, 1.050020598173771, 1.035890117219833, 1.023505699456694, 1.049895471432276, 1.022579227839277, 1.052410024451756, 1.011789568595402)
    "7" = @(1.02, 1.018525617602466, 1.047261280614948, 1.019701740947772, 1.049677367680689, 1.035844641014963, 1.023140497905607, 1.049716413478872, 1.02222692993797, 1.052126547647691, 1.011663939068427)
    "8" = @(1.02, 1.016512908571459, 1.046260974041784, 1.017968805612514, 1.048240391041728, 1.03565085203536, 1.021611002803895, 1.048964123001408, 1.020751679199596, 1.050938136683274, 1.011136899850583)
    "9" = @(1.02, 1.012956664017368, 1.0444887773233, 1.014908596520856, 1.045699330772224, 1.035295394798352, 1.01890409627993, 1.047623910025574, 1.018141547935889, 1.048830589368037, 1.010200797854806)
    "10" = @(1.02, 1.01057886840537, 1.043300874740713, 1.012863673752276, 1.043999150998451, 1.035049017934385, 1.017091233654087, 1.046720589361106, 1.016394024801468, 1.047416410908024, 1.009571643750741)
    "11" = @(1.02, 1.0095474893276, 1.042784964417849, 1.011976980884906, 1.04326146577713, 1.034940100732812, 1.016304203725232, 1.046327095892657, 1.015635487778754, 1.046801856978305, 1.009297980484599)
    "12" = @(1.02, 1.009164112130267, 1.042593100020983, 1.011647431885281, 1.042987228927309, 1.034899308052766, 1.016011550901997, 1.046180580854841, 1.015353449047237, 1.046573250277749, 1.009196141855619)
    "13" = @(1.02, 1.009246360540196, 1.042634266113093, 1.011718130061407, 1.04304606401449, 1.034908073429748, 1.016074340289043, 1.046212024865871, 1.015413960285702, 1.046622302349062, 1.009217995127632)
    "14" = @(1.02, 1.009515804959319, 1.042769109591558, 1.01194974421679, 1.043238801935954, 1.034936735652793, 1.016280019416042, 1.046314992142373, 1.015612180207779, 1.046782967106437, 1.009289566318406)
    "15" = @(1.02, 1.009681781542429, 1.042852160296442, 1.012092423708241, 1.043357523867945, 1.034954350864694, 1.016406703178162, 1.046378386734242, 1.015734272169935, 1.04688191362587, 1.009333638722408)
    "16" = @(1.02, 1.010647279329926, 1.043335081476528, 1.012922494145183, 1.044048076965312, 1.035056199309763, 1.017143422474447, 1.046746654634396, 1.016444326902632, 1.047457150136519, 1.009589779656025)
    "17" = @(1.02, 1.011252427571169, 1.043637592092244, 1.013442841526445, 1.044480840056042, 1.035119487736606, 1.017604993595315, 1.046977029565916, 1.016889226028765, 1.047817388490848, 1.009750117727743)
    "18" = @(1.02, 1.011605229439809, 1.04381389282556, 1.013746233382226, 1.044733119288688, 1.035156187239233, 1.017874023058829, 1.047111176729931, 1.017148550312981, 1.048027296830619, 1.009843521172261)
    "19" = @(1.02, 1.011725497079173, 1.04387398164128, 1.013849662375704, 1.044819115623741, 1.03516866426591, 1.017965721925431, 1.047156879035292, 1.017236943226875, 1.048098834169899, 1.009875349211112)
    "20" = @(1.02, 1.011187518637312, 1.043605150941936, 1.013387025432044, 1.044434423604186, 1.035112719786508, 1.01755549182724, 1.046952335965611, 1.016841511007119, 1.047778760334203, 1.00973292729965)
    "21" = @(1.02, 1.009436468023033, 1.042729407968616, 1.011881544974262, 1.043182051731683, 1.034928304619612, 1.01621946078874, 1.046284680627588, 1.01555381728333, 1.046735664581835, 1.009268495592347)
    "22" = @(1.02, 1.008333904875009, 1.042177448952804, 1.010933876359486, 1.042393316870428, 1.034810411446517, 1.015377619774331, 1.045862850301309, 1.014742544874372, 1.046077895059011, 1.008975400616499)
    "23" = @(1.02, 1.008918550242514, 1.042470180623124, 1.011436361409522, 1.042811566153862, 1.034873093219633, 1.015824071026052, 1.046086665103291, 1.015172774203533, 1.046426775337048, 1.009130879743878)
    "24" = @(1.02, 1.011216848718677, 1.04361981016091, 1.013412246687279, 1.044455397646697, 1.035115778598063, 1.017577860154867, 1.04696349464129, 1.016863071921271, 1.047796215392767, 1.009740695282076)
    "25" = @(1.02, 1.01387722827824, 1.044948065006041, 1.015700549318791, 1.046357327541651, 1.035388946115932, 1.019605325201447, 1.047972119356874, 1.018817614809866, 1.049377043632187, 1.010443690523745)
}

$colIndices = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

foreach ($rowKey in $newValues.Keys) {
    $rowNum = [int]$rowKey
    $rowVals = $newValues[$rowKey]
    for ($i = 0; $i -lt $colIndices.Length; $i++) {
        $ws.Cells.Item($rowNum, $colIndices[$i]).Value = $rowVals[$i]
    }
}
